# Update the "想去人数" (want-to-go count) figures (column F) across the
# "展览", "演出" and "全部类型" worksheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1294
$ws.Range("F4").Value  = 12998
$ws.Range("F13").Value = 4997
$ws.Range("F14").Value = 530
$ws.Range("F15").Value = 215
$ws.Range("F20").Value = 140
$ws.Range("F25").Value = 1316

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value  = 4454
$ws.Range("F11").Value = 366

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 1294
$ws.Range("F7").Value  = 12998
$ws.Range("F18").Value = 5000
$ws.Range("F19").Value = 530
$ws.Range("F20").Value = 4454
$ws.Range("F21").Value = 215
$ws.Range("F29").Value = 366
$ws.Range("F34").Value = 140
$ws.Range("F42").Value = 1316
